$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2032.9584
$ws.Range("J17").Value = 2069.1738
$ws.Range("L17").Value = 6207.5214
$ws.Range("N17").Value = -6543.5214
$ws.Range("H70").Value = 14400.75
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 18701
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 56103
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -56643
$ws.Range("H73").Value = 14400.75
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 18701
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 56103
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -57975
$ws.Range("H86").Value = 1783.2
$ws.Range("I86").Value = 1296
$ws.Range("J86").Value = 4950
$ws.Range("K86").Value = 1296
$ws.Range("L86").Value = 4950
$ws.Range("M86").Value = -173
$ws.Range("N86").Value = -7196
$ws.Range("H89").Value = 1783.2
$ws.Range("I89").Value = 1296
$ws.Range("J89").Value = 4950
$ws.Range("K89").Value = 6480
$ws.Range("L89").Value = 24750
$ws.Range("M89").Value = -864
$ws.Range("N89").Value = -35982
$ws.Range("H92").Value = 3521.0417
$ws.Range("I92").Value = 847.8946999999999
$ws.Range("K92").Value = 847.8946999999999
$ws.Range("M92").Value = 400.1053000000001
$ws.Range("H112").Value = 2464.1702
$ws.Range("J112").Value = 2245.561
$ws.Range("L112").Value = 6736.683000000001
$ws.Range("N112").Value = -8952.683000000001
$ws.Range("H113").Value = 2399.8333
$ws.Range("I113").Value = 2399.8333
$ws.Range("K113").Value = 2399.8333
$ws.Range("M113").Value = 854.1667000000002
$ws.Range("H116").Value = 5053968.5
$ws.Range("I116").Value = 7939647
$ws.Range("J116").Value = 4031.5
$ws.Range("K116").Value = 7939647
$ws.Range("L116").Value = 4031.5
$ws.Range("M116").Value = -7936205
$ws.Range("N116").Value = -10915.5
$ws.Range("H125").Value = 2083.923
$ws.Range("I125").Value = 1799.1428
$ws.Range("K125").Value = 16192.2852
$ws.Range("M125").Value = -13732.2852

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 10001749
$ws.Range("J11").Value = 375
$ws.Range("L11").Value = 375
$ws.Range("N11").Value = -663
$ws.Range("H32").Value = 5852803
$ws.Range("I32").Value = 3046.3572
$ws.Range("K32").Value = 3046.3572
$ws.Range("M32").Value = -2759.3572
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H61").Value = 4554.533
$ws.Range("I61").Value = 4200.579
$ws.Range("K61").Value = 4200.579
$ws.Range("M61").Value = -3988.579
$ws.Range("H110").Value = 4837.8237
$ws.Range("J110").Value = 6173.636
$ws.Range("L110").Value = 6173.636
$ws.Range("N110").Value = -10263.636
$ws.Range("H126").Value = 9990
$ws.Range("I126").Value = 9990
$ws.Range("K126").Value = 29970
$ws.Range("M126").Value = -27500
$ws.Range("H136").Value = 4554.533
$ws.Range("I136").Value = 4200.579
$ws.Range("K136").Value = 12601.737
$ws.Range("M136").Value = -10051.737

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 19235
$ws.Range("I26").Value = 19235
$ws.Range("K26").Value = 19235
$ws.Range("M26").Value = -18943
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H86").Value = 3186.8708
$ws.Range("I86").Value = 1565.5217
$ws.Range("J86").Value = 7848.25
$ws.Range("K86").Value = 1565.5217
$ws.Range("L86").Value = 7848.25
$ws.Range("M86").Value = -442.5217
$ws.Range("N86").Value = -10094.25
$ws.Range("H89").Value = 3186.8708
$ws.Range("I89").Value = 1565.5217
$ws.Range("J89").Value = 7848.25
$ws.Range("K89").Value = 7827.6085
$ws.Range("L89").Value = 39241.25
$ws.Range("M89").Value = -2211.6085
$ws.Range("N89").Value = -50473.25
$ws.Range("H123").Value = 80000
$ws.Range("J123").Value = 80000
$ws.Range("L123").Value = 80000
$ws.Range("N123").Value = -89800

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H99").Value = 6175417.5
$ws.Range("I99").Value = 11113461
$ws.Range("J99").Value = 2863.25
$ws.Range("K99").Value = 11113461
$ws.Range("L99").Value = 2863.25
$ws.Range("M99").Value = -11111963
$ws.Range("N99").Value = -5859.25
$ws.Range("H122").Value = 8931.462
$ws.Range("I122").Value = 4326.25
$ws.Range("J122").Value = 16299.8
$ws.Range("K122").Value = 12978.75
$ws.Range("L122").Value = 48899.39999999999
$ws.Range("M122").Value = -10528.75
$ws.Range("N122").Value = -53799.39999999999
$ws.Range("H126").Value = 6175417.5
$ws.Range("I126").Value = 11113461
$ws.Range("J126").Value = 2863.25
$ws.Range("K126").Value = 33340383
$ws.Range("L126").Value = 8589.75
$ws.Range("M126").Value = -33337913
$ws.Range("N126").Value = -13529.75
$ws.Range("H134").Value = 66675916
$ws.Range("I134").Value = 76927620
$ws.Range("K134").Value = 230782860
$ws.Range("M134").Value = -230780325

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I92").Value = 300
$ws.Range("K92").Value = 900
$ws.Range("M92").Value = 348
$ws.Range("H98").Value = 2188.1
$ws.Range("I98").Value = 731.1667
$ws.Range("J98").Value = 4373.5
$ws.Range("K98").Value = 2193.5001
$ws.Range("L98").Value = 13120.5
$ws.Range("M98").Value = -695.5001000000002
$ws.Range("N98").Value = -16116.5
$ws.Range("H109").Value = 13635.909
$ws.Range("J109").Value = 18333.334
$ws.Range("L109").Value = 55000.00199999999
$ws.Range("N109").Value = -57080.00199999999
$ws.Range("H136").Value = 55560856
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 3252.8333
$ws.Range("I137").Value = 1988.5
$ws.Range("K137").Value = 5965.5
$ws.Range("M137").Value = -865.5
$ws.Range("H139").Value = 16669125
$ws.Range("I139").Value = 20835132
$ws.Range("J139").Value = 5099
$ws.Range("K139").Value = 62505396
$ws.Range("L139").Value = 15297
$ws.Range("M139").Value = -62500256
$ws.Range("N139").Value = -25577

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4491.0293
$ws.Range("I80").Value = 3838
$ws.Range("J80").Value = 5071.5
$ws.Range("K80").Value = 3838
$ws.Range("L80").Value = 5071.5
$ws.Range("M80").Value = -2840
$ws.Range("N80").Value = -7067.5
$ws.Range("H83").Value = 4491.0293
$ws.Range("I83").Value = 3838
$ws.Range("J83").Value = 5071.5
$ws.Range("K83").Value = 19190
$ws.Range("L83").Value = 25357.5
$ws.Range("M83").Value = -14198
$ws.Range("N83").Value = -35341.5
$ws.Range("H133").Value = 90000
$ws.Range("J133").Value = 90000
$ws.Range("L133").Value = 90000
$ws.Range("N133").Value = -100120

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7943.778
$ws.Range("I7").Value = 5113.643
$ws.Range("K7").Value = 5113.643
$ws.Range("M7").Value = -5001.643
$ws.Range("H40").Value = 4577.8
$ws.Range("I40").Value = 4331.3335
$ws.Range("J40").Value = 4947.5
$ws.Range("K40").Value = 4331.3335
$ws.Range("L40").Value = 4947.5
$ws.Range("M40").Value = -4195.3335
$ws.Range("N40").Value = -5219.5
$ws.Range("H41").Value = 29749.5
$ws.Range("I41").Value = 49499
$ws.Range("J41").Value = 10000
$ws.Range("K41").Value = 49499
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = -49061
$ws.Range("N41").Value = -10876
$ws.Range("H50").Value = 37827.332
$ws.Range("I50").Value = 40749
$ws.Range("J50").Value = 31984
$ws.Range("K50").Value = 40749
$ws.Range("L50").Value = 31984
$ws.Range("M50").Value = -40112
$ws.Range("N50").Value = -33258
$ws.Range("H123").Value = 10000
$ws.Range("J123").Value = 10000
$ws.Range("L123").Value = 10000
$ws.Range("N123").Value = -19800
$ws.Range("H126").Value = 7943.778
$ws.Range("I126").Value = 5113.643
$ws.Range("K126").Value = 15340.929
$ws.Range("M126").Value = -12870.929

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 30000
$ws.Range("I34").Value = 38000
$ws.Range("K34").Value = 38000
$ws.Range("M34").Value = -37797
$ws.Range("H37").Value = 19999.5
$ws.Range("I37").Value = 14999
$ws.Range("K37").Value = 14999
$ws.Range("M37").Value = -14796
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H81").Value = 1627
$ws.Range("I81").Value = 1523.1666
$ws.Range("K81").Value = 3046.3332
$ws.Range("M81").Value = -1985.3332
$ws.Range("H84").Value = 1627
$ws.Range("I84").Value = 1523.1666
$ws.Range("K84").Value = 15231.666
$ws.Range("M84").Value = -9927.666000000001
$ws.Range("H122").Value = 13903.074
$ws.Range("J122").Value = 20144.834
$ws.Range("L122").Value = 60434.50199999999
$ws.Range("N122").Value = -65334.50199999999
$ws.Range("H132").Value = 7787.2144
$ws.Range("I132").Value = 5813.52
$ws.Range("K132").Value = 17440.56
$ws.Range("M132").Value = -14910.56
